# Atualização de bases das ligas, do dia: 10-06-2024 às 21:53
# Fix two swapped pairs of match rows: the rows' "id" sequence numbers
# (column A) were correct, but the actual match data (columns B:AD) had
# been written into the wrong row of each pair. Swap the B:AD payloads
# back between each pair of rows so the data lines up with the right
# sequence number again.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$pairs = @(
    @(128, 129),
    @(287, 288),
    @(296, 297)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $rng1 = $ws.Range("B$r1`:AD$r1")
    $rng2 = $ws.Range("B$r2`:AD$r2")

    $data1 = $rng1.Value2
    $data2 = $rng2.Value2

    $rng1.Value2 = $data2
    $rng2.Value2 = $data1
}
